$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.030.28"
$ws.Range("E2").Value = "  +0.35%  "
$ws.Range("D3").Value = "2.561.03"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'582.79"
$ws.Range("E5").Value = "  +0.82%  "
$ws.Range("D6").Value = "'171.40"
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.520"
$ws.Range("E8").Value = "  +1.70%  "
$ws.Range("D9").Value = "2.559.16"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("E12").Value = "  +3.16%  "
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("D14").Value = "3.032.29"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "71.129.62"
$ws.Range("E15").Value = "  +0.66%  "
$ws.Range("E16").Value = "  -2.60%  "
$ws.Range("D17").Value = "'25.53"
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("D18").Value = "2.566.40"
$ws.Range("E18").Value = "  -0.59%  "
$ws.Range("D19").Value = "'11.65"
$ws.Range("E19").Value = "  -1.22%  "
$ws.Range("E20").Value = "  +3.85%  "
$ws.Range("D21").Value = "'357.24"
$ws.Range("E21").Value = "  -1.92%  "
$ws.Range("E22").Value = "  -0.83%  "
$ws.Range("D23").Value = "'2.06"
$ws.Range("E23").Value = "  +3.79%  "
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").Value = "'70.38"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").Value = "'4.10"
$ws.Range("E26").Value = "  -1.24%  "
$ws.Range("E27").Value = "  -1.42%  "
$ws.Range("D28").Value = "2.706.93"
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("D30").Value = "0.0₃0926"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("D31").Value = "'8.00"
$ws.Range("E31").Value = "  +2.65%  "
$ws.Range("D32").Value = "'474.13"
$ws.Range("E32").Value = "  -1.94%  "
$ws.Range("D33").Value = "'1.28"
$ws.Range("E33").Value = "  -1.21%  "
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").Value = "'0.118"
$ws.Range("E36").Value = "  +3.75%  "
$ws.Range("D37").Value = "'157.26"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "'18.88"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("D39").Value = "'19.12"
$ws.Range("E39").Value = "  +1.46%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("D42").Value = "'0.322"
$ws.Range("E42").Value = "  +0.81%  "
$ws.Range("E43").Value = "  -3.95%  "
$ws.Range("E44").Value = "  -3.83%  "
$ws.Range("E45").Value = "  -11.59%  "
$ws.Range("D46").Value = "'38.76"
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("D47").Value = "'145.92"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").Value = "'0.540"
$ws.Range("E48").Value = "  +1.79%  "
$ws.Range("E49").Value = "  -0.01%  "
$ws.Range("D50").Value = "'1.62"
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("D51").Value = "'0.0741"
$ws.Range("E51").Value = "  +1.20%  "
